$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-07 Tuesday" "2025-10-08 Wednesday"

Replace-Text "74×33=2442" "56×54=3024"
Replace-Text "94×58=5452" "52×16=832"
Replace-Text "85×65=5525" "15×77=1155"
Replace-Text "80×80=6400" "69×42=2898"
Replace-Text "14×35=490" "96×42=4032"

Replace-Text "57×37=2109" "12×23=276"
Replace-Text "23×43=989" "31×20=620"
Replace-Text "61×49=2989" "61×81=4941"
Replace-Text "63×14=882" "25×49=1225"
Replace-Text "63×20=1260" "16×82=1312"

Replace-Text "21×28=588" "33×30=990"
Replace-Text "61×71=4331" "11×51=561"
Replace-Text "89×17=1513" "47×73=3431"
Replace-Text "47×53=2491" "57×55=3135"
Replace-Text "48×72=3456" "28×56=1568"

Replace-Text "52×69=3588" "60×12=720"
Replace-Text "66×22=1452" "59×12=708"
Replace-Text "56×80=4480" "42×41=1722"
Replace-Text "40×23=920" "83×20=1660"
Replace-Text "41×79=3239" "70×70=4900"

Replace-Text "64×69=4416" "84×11=924"
Replace-Text "20×13=260" "46×93=4278"
Replace-Text "63×41=2583" "55×62=3410"
Replace-Text "69×96=6624" "98×86=8428"
Replace-Text "59×99=5841" "96×90=8640"
